$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = "sv"
$ws.Range("J2").Value = "Statement-opinion"
$ws.Range("I3").Value = "b"
$ws.Range("J3").Value = "Acknowledge (Backchannel)"
$ws.Range("I7").Value = "ba"
$ws.Range("J7").Value = "Appreciation"
$ws.Range("I9").Value = "sv"
$ws.Range("J9").Value = "Statement-opinion"
$ws.Range("I14").Value = "sd"
$ws.Range("J14").Value = "Statement-non-opinion"
$ws.Range("I17").Value = "b"
$ws.Range("J17").Value = "Acknowledge (Backchannel)"
$ws.Range("I27").Value = "ba"
$ws.Range("J27").Value = "Appreciation"
$ws.Range("I32").Value = "ba"
$ws.Range("J32").Value = "Appreciation"
$ws.Range("I35").Value = "ba"
$ws.Range("J35").Value = "Appreciation"
$ws.Range("I38").Value = "ba"
$ws.Range("J38").Value = "Appreciation"
$ws.Range("I41").Value = "sv"
$ws.Range("J41").Value = "Statement-opinion"
$ws.Range("I51").Value = "sv"
$ws.Range("J51").Value = "Statement-opinion"
$ws.Range("I72").Value = "ba"
$ws.Range("J72").Value = "Appreciation"
$ws.Range("I75").Value = "aa"
$ws.Range("J75").Value = "Agree/Accept"
$ws.Range("I76").Value = "ba"
$ws.Range("J76").Value = "Appreciation"
$ws.Range("I81").Value = "sd"
$ws.Range("J81").Value = "Statement-non-opinion"
$ws.Range("I84").Value = "%"
$ws.Range("J84").Value = "Uninterpretable"
$ws.Range("I87").Value = "aa"
$ws.Range("J87").Value = "Agree/Accept"
$ws.Range("I89").Value = "ba"
$ws.Range("J89").Value = "Appreciation"
$ws.Range("I102").Value = "ba"
$ws.Range("J102").Value = "Appreciation"
$ws.Range("I107").Value = "ba"
$ws.Range("J107").Value = "Appreciation"
$ws.Range("I109").Value = "sv"
$ws.Range("J109").Value = "Statement-opinion"
$ws.Range("I117").Value = "sd"
$ws.Range("J117").Value = "Statement-non-opinion"
$ws.Range("I124").Value = "b"
$ws.Range("J124").Value = "Acknowledge (Backchannel)"
$ws.Range("I126").Value = "sv"
$ws.Range("J126").Value = "Statement-opinion"
$ws.Range("I127").Value = "sv"
$ws.Range("J127").Value = "Statement-opinion"
$ws.Range("I139").Value = "ba"
$ws.Range("J139").Value = "Appreciation"
$ws.Range("I140").Value = "sv"
$ws.Range("J140").Value = "Statement-opinion"
$ws.Range("I144").Value = "aa"
$ws.Range("J144").Value = "Agree/Accept"
$ws.Range("I152").Value = "b"
$ws.Range("J152").Value = "Acknowledge (Backchannel)"
$ws.Range("I166").Value = "ba"
$ws.Range("J166").Value = "Appreciation"
$ws.Range("I174").Value = "b"
$ws.Range("J174").Value = "Acknowledge (Backchannel)"
$ws.Range("I188").Value = "ba"
$ws.Range("J188").Value = "Appreciation"
$ws.Range("I193").Value = "aa"
$ws.Range("J193").Value = "Agree/Accept"
$ws.Range("I200").Value = "%"
$ws.Range("J200").Value = "Uninterpretable"
$ws.Range("I205").Value = "sd"
$ws.Range("J205").Value = "Statement-non-opinion"
$ws.Range("I209").Value = "b"
$ws.Range("J209").Value = "Acknowledge (Backchannel)"
$ws.Range("I218").Value = "ba"
$ws.Range("J218").Value = "Appreciation"
$ws.Range("I246").Value = "b"
$ws.Range("J246").Value = "Acknowledge (Backchannel)"
$ws.Range("I272").Value = "ba"
$ws.Range("J272").Value = "Appreciation"
$ws.Range("I275").Value = "b"
$ws.Range("J275").Value = "Acknowledge (Backchannel)"
$ws.Range("I276").Value = "b"
$ws.Range("J276").Value = "Acknowledge (Backchannel)"
$ws.Range("I277").Value = "b"
$ws.Range("J277").Value = "Acknowledge (Backchannel)"
$ws.Range("I279").Value = "sd"
$ws.Range("J279").Value = "Statement-non-opinion"
$ws.Range("I295").Value = "sv"
$ws.Range("J295").Value = "Statement-opinion"
$ws.Range("I306").Value = "b"
$ws.Range("J306").Value = "Acknowledge (Backchannel)"
$ws.Range("I317").Value = "sv"
$ws.Range("J317").Value = "Statement-opinion"
$ws.Range("I327").Value = "ba"
$ws.Range("J327").Value = "Appreciation"
$ws.Range("I331").Value = "b"
$ws.Range("J331").Value = "Acknowledge (Backchannel)"
$ws.Range("I333").Value = "sd"
$ws.Range("J333").Value = "Statement-non-opinion"
$ws.Range("I340").Value = "aa"
$ws.Range("J340").Value = "Agree/Accept"
$ws.Range("I341").Value = "sv"
$ws.Range("J341").Value = "Statement-opinion"
$ws.Range("I343").Value = "b"
$ws.Range("J343").Value = "Acknowledge (Backchannel)"
$ws.Range("I353").Value = "ba"
$ws.Range("J353").Value = "Appreciation"
$ws.Range("I354").Value = "sd"
$ws.Range("J354").Value = "Statement-non-opinion"
$ws.Range("I357").Value = "%"
$ws.Range("J357").Value = "Uninterpretable"
$ws.Range("I358").Value = "sd"
$ws.Range("J358").Value = "Statement-non-opinion"
$ws.Range("I360").Value = "sd"
$ws.Range("J360").Value = "Statement-non-opinion"
$ws.Range("I363").Value = "sv"
$ws.Range("J363").Value = "Statement-opinion"
$ws.Range("I370").Value = "%"
$ws.Range("J370").Value = "Uninterpretable"
$ws.Range("I373").Value = "sd"
$ws.Range("J373").Value = "Statement-non-opinion"
$ws.Range("I376").Value = "b"
$ws.Range("J376").Value = "Acknowledge (Backchannel)"
$ws.Range("I379").Value = "aa"
$ws.Range("J379").Value = "Agree/Accept"
$ws.Range("I381").Value = "aa"
$ws.Range("J381").Value = "Agree/Accept"
$ws.Range("I385").Value = "sd"
$ws.Range("J385").Value = "Statement-non-opinion"
$ws.Range("I398").Value = "ba"
$ws.Range("J398").Value = "Appreciation"
$ws.Range("I407").Value = "b"
$ws.Range("J407").Value = "Acknowledge (Backchannel)"
$ws.Range("I409").Value = "aa"
$ws.Range("J409").Value = "Agree/Accept"
$ws.Range("I411").Value = "b"
$ws.Range("J411").Value = "Acknowledge (Backchannel)"
$ws.Range("I413").Value = "b"
$ws.Range("J413").Value = "Acknowledge (Backchannel)"
$ws.Range("I415").Value = "ba"
$ws.Range("J415").Value = "Appreciation"
$ws.Range("I416").Value = "ba"
$ws.Range("J416").Value = "Appreciation"
$ws.Range("I418").Value = "sd"
$ws.Range("J418").Value = "Statement-non-opinion"
$ws.Range("I421").Value = "aa"
$ws.Range("J421").Value = "Agree/Accept"
$ws.Range("I424").Value = "b"
$ws.Range("J424").Value = "Acknowledge (Backchannel)"
$ws.Range("I427").Value = "aa"
$ws.Range("J427").Value = "Agree/Accept"
$ws.Range("I428").Value = "sd"
$ws.Range("J428").Value = "Statement-non-opinion"
$ws.Range("I439").Value = "aa"
$ws.Range("J439").Value = "Agree/Accept"
$ws.Range("I451").Value = "ba"
$ws.Range("J451").Value = "Appreciation"
$ws.Range("I453").Value = "ba"
$ws.Range("J453").Value = "Appreciation"
$ws.Range("I455").Value = "b"
$ws.Range("J455").Value = "Acknowledge (Backchannel)"
$ws.Range("I466").Value = "qy"
$ws.Range("J466").Value = "Yes-No-Question"